{"js": "// Update the date line and the twenty-five \"A\u00d7B=\" practice problems.\n// Each old value is unique in the document, so a simple search +\n// replace (preserving the run's existing formatting) for every pair\n// is sufficient and avoids any row/column index bookkeeping.\nconst replacements = [\n  [\"2025-10-14 Tuesday\", \"2025-10-15 Wednesday\"],\n  [\"987\u00d74=\", \"111\u00d75=\"],\n  [\"232\u00d77=\", \"841\u00d78=\"],\n  [\"645\u00d77=\", \"680\u00d75=\"],\n  [\"498\u00d72=\", \"274\u00d77=\"],\n  [\"306\u00d75=\", \"942\u00d76=\"],\n  [\"831\u00d76=\", \"149\u00d78=\"],\n  [\"386\u00d76=\", \"770\u00d78=\"],\n  [\"350\u00d74=\", \"619\u00d75=\"],\n  [\"872\u00d73=\", \"872\u00d79=\"],\n  [\"742\u00d78=\", \"579\u00d76=\"],\n  [\"892\u00d73=\", \"855\u00d73=\"],\n  [\"891\u00d73=\", \"936\u00d73=\"],\n  [\"676\u00d75=\", \"732\u00d78=\"],\n  [\"121\u00d75=\", \"696\u00d76=\"],\n  [\"832\u00d73=\", \"612\u00d78=\"],\n  [\"566\u00d79=\", \"969\u00d72=\"],\n  [\"894\u00d72=\", \"662\u00d73=\"],\n  [\"857\u00d77=\", \"860\u00d75=\"],\n  [\"121\u00d76=\", \"698\u00d76=\"],\n  [\"784\u00d75=\", \"672\u00d79=\"],\n  [\"634\u00d75=\", \"925\u00d75=\"],\n  [\"998\u00d74=\", \"178\u00d72=\"],\n  [\"922\u00d77=\", \"335\u00d78=\"],\n  [\"233\u00d78=\", \"721\u00d75=\"],\n  [\"900\u00d74=\", \"155\u00d74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the twenty-five \"A\u00d7B=\" practice problems.\n# Every old value is unique in the document, so Find/Replace across the\n# whole document body for each pair is sufficient and preserves the\n# existing run formatting (font, size) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-14 Tuesday\", \"2025-10-15 Wednesday\"),\n    @(\"987\u00d74=\", \"111\u00d75=\"),\n    @(\"232\u00d77=\", \"841\u00d78=\"),\n    @(\"645\u00d77=\", \"680\u00d75=\"),\n    @(\"498\u00d72=\", \"274\u00d77=\"),\n    @(\"306\u00d75=\", \"942\u00d76=\"),\n    @(\"831\u00d76=\", \"149\u00d78=\"),\n    @(\"386\u00d76=\", \"770\u00d78=\"),\n    @(\"350\u00d74=\", \"619\u00d75=\"),\n    @(\"872\u00d73=\", \"872\u00d79=\"),\n    @(\"742\u00d78=\", \"579\u00d76=\"),\n    @(\"892\u00d73=\", \"855\u00d73=\"),\n    @(\"891\u00d73=\", \"936\u00d73=\"),\n    @(\"676\u00d75=\", \"732\u00d78=\"),\n    @(\"121\u00d75=\", \"696\u00d76=\"),\n    @(\"832\u00d73=\", \"612\u00d78=\"),\n    @(\"566\u00d79=\", \"969\u00d72=\"),\n    @(\"894\u00d72=\", \"662\u00d73=\"),\n    @(\"857\u00d77=\", \"860\u00d75=\"),\n    @(\"121\u00d76=\", \"698\u00d76=\"),\n    @(\"784\u00d75=\", \"672\u00d79=\"),\n    @(\"634\u00d75=\", \"925\u00d75=\"),\n    @(\"998\u00d74=\", \"178\u00d72=\"),\n    @(\"922\u00d77=\", \"335\u00d78=\"),\n    @(\"233\u00d78=\", \"721\u00d75=\"),\n    @(\"900\u00d74=\", \"155\u00d74=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n}\n"}
